$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped crypto data.
# D-column price values are forced to Text format since they are pre-formatted
# strings (e.g. "36.90", "0.950") that must keep their exact trailing zeros/
# grouping dots rather than being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.014.93'
$ws.Range("E2").Value = '  -2.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.908.14'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '365.01'
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.11'
$ws.Range("E6").Value = '  -6.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.539'
$ws.Range("E7").Value = '  -5.18%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -7.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.90'
$ws.Range("E10").Value = '  -5.75%  '
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("E12").Value = '  -4.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.37'
$ws.Range("E13").Value = '  -6.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.366.36'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.35'
$ws.Range("E15").Value = '  -6.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.901.30'
$ws.Range("E16").Value = '  -1.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.950'
$ws.Range("E17").Value = '  -3.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.999.53'
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.29'
$ws.Range("E19").Value = '  -6.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.24'
$ws.Range("E20").Value = '  -4.56%  '
$ws.Range("E21").Value = '  -7.20%  '
$ws.Range("E22").Value = '  -4.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.07'
$ws.Range("E23").Value = '  -3.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.12'
$ws.Range("E24").Value = '  -3.37%  '
$ws.Range("E25").Value = '  -4.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.35'
$ws.Range("E26").Value = '  +4.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.174'
$ws.Range("E27").Value = '  -5.84%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.89'
$ws.Range("E29").Value = '  -4.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.29'
$ws.Range("E30").Value = '  -6.81%  '
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.91'
$ws.Range("E32").Value = '  -5.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.09'
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("E34").Value = '  -6.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '34.93'
$ws.Range("E35").Value = '  -7.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.58'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E38").Value = '  -5.30%  '
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.13'
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.88'
$ws.Range("E41").Value = '  -7.91%  '
$ws.Range("E42").Value = '  -7.33%  '
$ws.Range("E43").Value = '  -5.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.48'
$ws.Range("E44").Value = '  -3.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '117.73'
$ws.Range("E45").Value = '  -1.69%  '
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.062.55'
$ws.Range("E47").Value = '  -3.30%  '
$ws.Range("E48").Value = '  -7.99%  '
$ws.Range("E49").Value = '  -8.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.200.56'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.234'
$ws.Range("E51").Value = '  -7.15%  '
